$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MuniEntryPleas")

# Insert a new column before D (so old D/E shift to E/F)
$ws.Columns("D").Insert()

# Header for new column
$ws.Range("D1").Value = "CaseTypeCode"

# Map of case-number prefix codes per row (derived from CaseNumber in column A)
$codes = @{2="TRD"; 3="TRC"; 4="CRB"; 5="CRB"; 6="TRC"; 7="TRC"; 8="TRC"; 9="TRC"; 10="TRD"; 11="CRB"; 12="TRC"}

foreach ($r in 2..12) {
    $ws.Range("D$r").Value = $codes[$r]
}

# Update selection as in the target file
$ws.Range("D13").Select()
